# Update profit files after running on 2025-12-24
# Appends the new daily row (row 30) to the profit data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (matching the existing rows,
# which store dates as plain strings rather than Excel date serials).
# A leading apostrophe forces text entry instead of date auto-detection;
# ClearFormats() then drops the "quote prefix" cell style that entering
# text this way would otherwise leave behind, so the new cell ends up
# with the same (default) styling as the rest of the data rows.
$ws.Range("A30").Value = "'12/24/2025"
$ws.Range("A30").ClearFormats()

$ws.Range("B30").Value = 12129.23
$ws.Range("C30").Value = 0.2070483707744404
$ws.Range("D30").Value = 0.7929516292255596
$ws.Range("E30").Value = -136.51
$ws.Range("F30").Value = -28.15
$ws.Range("G30").Value = -20923.87
$ws.Range("H30").Value = -68.51000000000001
$ws.Range("I30").Value = -491.51
$ws.Range("J30").Value = -16.37
